$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("O2").Value = 210210930.23
$ws.Range("P2").Value = 44315188.57
$ws.Range("Q2").Value = 86032857.31
$ws.Range("R2").Value = 59.2168693808
$ws.Range("S2").Value = 58223610.26
$ws.Range("T2").Value = 3.8537564093
$ws.Range("U2").Value = 1938497.61
$ws.Range("V2").Value = -87.8719365387
$ws.Range("W2").Value = 95905289.41
$ws.Range("X2").Value = 39690240.28
$ws.Range("Y2").Value = -32.8873378024
$ws.Range("Z2").Value = 7078071.53
$ws.Range("AA2").Value = -63.7000751326
$ws.Range("AB2").Value = 114305640.82
$ws.Range("AC2").Value = 48.0764859734
$ws.Range("AD2").Value = 12.5684571318
$ws.Range("AE2").Value = -12.4527418248
$ws.Range("AF2").Value = 172.3477108213
$ws.Range("AG2").Value = 45.6233599771
